$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.0066958374324755, 1.241785817101409, 0.3137894063765676, 3.1100811202262877)"
$ws.Range("C2").Value = "JSU(-1.4735914259337846, 1.2046027778305293, 1.8890185999045197, 4.773950464883239)"
$ws.Range("D2").Value = "JSU(-0.7685359356992993, 1.0473209131217454, 0.7831275714411181, 2.59247617602525)"
$ws.Range("E2").Value = "NIG(1.7211183075137435, 1.3220733869818715, 3.092571613984099, 6.512405539854512)"
